$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.835941000000001
$ws.Range("H2").Value = 17.507823
$ws.Range("I2").Value = 0.03643643319117328
$ws.Range("J2").Value = 0.03643643319117327
$ws.Range("M2").Value = 1.847479
$ws.Range("N2").Value = 5.542437
$ws.Range("O2").Value = 0.3181373042830636
$ws.Range("P2").Value = 0.3181373042830637
$ws.Range("Q2").Value = 10.781778442739
$ws.Range("R2").Value = 97.036005984651
$ws.Range("S2").Value = 0.01159178863312981
$ws.Range("T2").Value = 0.01159178863312981

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.835941000000001
$ws.Range("H3").Value = 17.507823
$ws.Range("I3").Value = 0.03643643319117328
$ws.Range("J3").Value = 0.03643643319117327
$ws.Range("N3").Value = 6.300930999999999
$ws.Range("O3").Value = 0.3616750542791174
$ws.Range("P3").Value = 0.3616750542791174
$ws.Range("Q3").Value = 12.25728718702367
$ws.Range("R3").Value = 110.315584683213
$ws.Range("S3").Value = 0.01317814895215503
$ws.Range("T3").Value = 0.01317814895215503

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.835941000000001
$ws.Range("H4").Value = 17.507823
$ws.Range("I4").Value = 0.03643643319117328
$ws.Range("J4").Value = 0.03643643319117327
$ws.Range("M4").Value = 1.859385666666666
$ws.Range("N4").Value = 5.578156999999999
$ws.Range("O4").Value = 0.320187641437819
$ws.Range("P4").Value = 0.320187641437819
$ws.Range("Q4").Value = 10.85126504691233
$ws.Range("R4").Value = 97.661385422211
$ws.Range("S4").Value = 0.01166649560588844
$ws.Range("T4").Value = 0.01166649560588843

# Row 5
$ws.Range("G5").Value = 17.50798033333334
$ws.Range("H5").Value = 52.52394100000001
$ws.Range("I5").Value = 0.1093102818770573
$ws.Range("J5").Value = 0.1093102818770573
$ws.Range("M5").Value = 1.847479
$ws.Range("N5").Value = 5.542437
$ws.Range("O5").Value = 0.3181373042830636
$ws.Range("P5").Value = 0.3181373042830637
$ws.Range("Q5").Value = 32.34562599824633
$ws.Range("R5").Value = 291.110633984217
$ws.Range("S5").Value = 0.03477567840678884
$ws.Range("T5").Value = 0.03477567840678884

# Row 6
$ws.Range("G6").Value = 17.50798033333334
$ws.Range("H6").Value = 52.52394100000001
$ws.Range("I6").Value = 0.1093102818770573
$ws.Range("J6").Value = 0.1093102818770573
$ws.Range("N6").Value = 6.300930999999999
$ws.Range("O6").Value = 0.3616750542791174
$ws.Range("P6").Value = 0.3616750542791174
$ws.Range("S6").Value = 0.03953480213115032
$ws.Range("T6").Value = 0.03953480213115031

# Row 7
$ws.Range("G7").Value = 17.50798033333334
$ws.Range("H7").Value = 52.52394100000001
$ws.Range("I7").Value = 0.1093102818770573
$ws.Range("J7").Value = 0.1093102818770573
$ws.Range("M7").Value = 1.859385666666666
$ws.Range("N7").Value = 5.578156999999999
$ws.Range("O7").Value = 0.320187641437819
$ws.Range("P7").Value = 0.320187641437819
$ws.Range("Q7").Value = 32.55408768408189
$ws.Range("R7").Value = 292.986789156737
$ws.Range("S7").Value = 0.03499980133911815
$ws.Range("T7").Value = 0.03499980133911815

# Row 8
$ws.Range("G8").Value = 136.8238143333333
$ws.Range("H8").Value = 410.471443
$ws.Range("I8").Value = 0.8542532849317694
$ws.Range("J8").Value = 0.8542532849317694
$ws.Range("M8").Value = 1.847479
$ws.Range("N8").Value = 5.542437
$ws.Range("O8").Value = 0.3181373042830636
$ws.Range("P8").Value = 0.3181373042830637
$ws.Range("Q8").Value = 252.7791236807323
$ws.Range("R8").Value = 2275.012113126591
$ws.Range("S8").Value = 0.271769837243145
$ws.Range("T8").Value = 0.271769837243145

# Row 9
$ws.Range("G9").Value = 136.8238143333333
$ws.Range("H9").Value = 410.471443
$ws.Range("I9").Value = 0.8542532849317694
$ws.Range("J9").Value = 0.8542532849317694
$ws.Range("N9").Value = 6.300930999999999
$ws.Range("O9").Value = 0.3616750542791174
$ws.Range("P9").Value = 0.3616750542791174
$ws.Range("Q9").Value = 287.3724710903815
$ws.Range("R9").Value = 2586.352239813433
$ws.Range("S9").Value = 0.308962103195812
$ws.Range("T9").Value = 0.308962103195812

# Row 10
$ws.Range("G10").Value = 136.8238143333333
$ws.Range("H10").Value = 410.471443
$ws.Range("I10").Value = 0.8542532849317694
$ws.Range("J10").Value = 0.8542532849317694
$ws.Range("M10").Value = 1.859385666666666
$ws.Range("N10").Value = 5.578156999999999
$ws.Range("O10").Value = 0.320187641437819
$ws.Range("P10").Value = 0.320187641437819
$ws.Range("Q10").Value = 254.4082392300612
$ws.Range("R10").Value = 2289.674153070551
$ws.Range("S10").Value = 0.2735213444928125
$ws.Range("T10").Value = 0.2735213444928125

